$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 3-8 (columns D, M, N, O, P, Q, S) are being permuted in a
# 6-cycle: new row r receives the old values that used to live in row
# $map[r]. Capture the "before" snapshot first so the in-place writes
# below don't clobber source data we still need to read.

$map = @{3=5; 4=7; 5=8; 6=3; 7=6; 8=4}

$snapshot = @{}
foreach ($r in 3..8) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

foreach ($r in 3..8) {
    $src = $snapshot[$map[$r]]
    $ws.Cells.Item($r, 4).Value2 = $src.D
    $ws.Cells.Item($r, 13).Value2 = $src.M
    $ws.Cells.Item($r, 14).Value2 = $src.N
    $ws.Cells.Item($r, 15).Value2 = $src.O
    $ws.Cells.Item($r, 16).Value2 = $src.P
    $ws.Cells.Item($r, 17).Value2 = $src.Q
    $ws.Cells.Item($r, 19).Value2 = $src.S
}
